$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 25: new purchase recorded on 05/17/2025 (date stored as plain text,
# matching the convention used by the other recently-appended rows).
$ws.Cells.Item(25, 1).Value = "'05/17/2025"
$ws.Cells.Item(25, 1).Style = "Normal"
$ws.Cells.Item(25, 2).Value = 448.8969999999972
$ws.Cells.Item(25, 3).Value = 0.1113841259799025
$ws.Cells.Item(25, 4).Value = 50
